# Actualización automática 2025-08-25 12:40:09
#
# Three asesor/cliente rows are being removed from the "VENTAS POR GRUPO"
# and "VENTA MENSUAL" sheets:
#   - ASES GAVILANEZ FAUSTO HERNAN        (row 9)
#   - BRAVO MONTENEGRO DANIEL ANDRES      (row 12)
#   - LATACELA ZUÑIGA JUAN FERNANDO       (row 28)
#
# Removing these rows shifts everything below them up, shrinking the used
# range from row 57 down to row 54 on both sheets. The trailing summary
# row (counts "N de 55" on sheet1, totals on sheet2) needs to be updated
# by hand afterwards because those cells hold cached literal values, not
# live formulas.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Delete bottom-most row first so the remaining row numbers (9 and 12)
# still point at the rows we intend to remove.
$ws1.Rows.Item(28).Delete()
$ws1.Rows.Item(12).Delete()
$ws1.Rows.Item(9).Delete()

$ws2.Rows.Item(28).Delete()
$ws2.Rows.Item(12).Delete()
$ws2.Rows.Item(9).Delete()

# --- Fix up the trailing summary row on "VENTAS POR GRUPO" (now row 54) ---
# The cells held text like "3 de 55"; the denominator must become 52
# (55 rows of data - 3 deleted = 52) while the numerators stay the same.
$sheet1Counts = @{
    "C" = 3
    "D" = 10
    "E" = 4
    "F" = 0
    "G" = 0
    "H" = 3
    "I" = 4
    "J" = 0
    "K" = 1
    "L" = 5
    "M" = 14
    "N" = 0
    "O" = 5
    "P" = 3
    "Q" = 0
    "R" = 1
}
foreach ($col in $sheet1Counts.Keys) {
    $n = $sheet1Counts[$col]
    $ws1.Range("${col}54").Value = "$n de 52"
}

# --- Fix up the trailing totals row on "VENTA MENSUAL" (now row 54) ---
$sheet2Totals = @{
    "C" = 97519.63
    "D" = 88077.29000000001
    "E" = 93548.25
    "F" = 62018.24
    "G" = 128670.11
}
foreach ($col in $sheet2Totals.Keys) {
    $ws2.Range("${col}54").Value = $sheet2Totals[$col]
}
